$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2: client name changes from "Exemplo" to "Cliente A"
$ws.Range("A2").Value = "Cliente A"
$ws.Range("B2").Value = "Ligar"
$ws.Range("C2").Value = "A Fazer"

# Add new row 3 for the second client
$ws.Range("A3").Value = "Cliente B"
$ws.Range("B3").Value = "Enviar e-mail"
$ws.Range("C3").Value = "Em andamento"
